$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '244.79'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '-0.95%'
$c.Style = 'Normal'
$c = $ws.Range('G2')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '4.50%'
$c.Style = 'Normal'
$c = $ws.Range('G3')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '5.115'
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '0.53%'
$c.Style = 'Normal'
$c = $ws.Range('G4')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '0.05687'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '1.50%'
$c.Style = 'Normal'
$c = $ws.Range('G5')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '0.66%'
$c.Style = 'Normal'
$c = $ws.Range('G6')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.8204'
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '0.88%'
$c.Style = 'Normal'
$c = $ws.Range('G7')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.8543'
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '1.23%'
$c.Style = 'Normal'
$c = $ws.Range('G8')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B9')
$c.NumberFormat = '@'
$c.Value = 'WazirX'
$c.Style = 'Normal'
$c = $ws.Range('C9')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.1334'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '0.06%'
$c.Style = 'Normal'
$c = $ws.Range('G9')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B10')
$c.NumberFormat = '@'
$c.Value = 'MandalaExchangeToken'
$c.Style = 'Normal'
$c = $ws.Range('C10')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.06940'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '-0.64%'
$c.Style = 'Normal'
$c = $ws.Range('G10')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B11')
$c.NumberFormat = '@'
$c.Value = 'LiechtensteinCryptoassetsExchange'
$c.Style = 'Normal'
$c = $ws.Range('C11')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.03216'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '0.31%'
$c.Style = 'Normal'
$c = $ws.Range('G11')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '1.16%'
$c.Style = 'Normal'
$c = $ws.Range('G12')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.09396'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '0.15%'
$c.Style = 'Normal'
$c = $ws.Range('G13')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.001517'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '0.27%'
$c.Style = 'Normal'
$c = $ws.Range('G14')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.04010'
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '-13.67%'
$c.Style = 'Normal'
$c = $ws.Range('G15')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B16')
$c.NumberFormat = '@'
$c.Value = 'One'
$c.Style = 'Normal'
$c = $ws.Range('C16')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0006028'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '-93.95%'
$c.Style = 'Normal'
$c = $ws.Range('G16')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B17')
$c.NumberFormat = '@'
$c.Value = 'TigerCash'
$c.Style = 'Normal'
$c = $ws.Range('C17')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.006213'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '1.43%'
$c.Style = 'Normal'
$c = $ws.Range('G17')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B18')
$c.NumberFormat = '@'
$c.Value = 'LEO'
$c.Style = 'Normal'
$c = $ws.Range('C18')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '3.512'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '-2.68%'
$c.Style = 'Normal'
$c = $ws.Range('G18')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B19')
$c.NumberFormat = '@'
$c.Value = 'GateToken'
$c.Style = 'Normal'
$c = $ws.Range('C19')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '3.010'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '-0.36%'
$c.Style = 'Normal'
$c = $ws.Range('G19')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B20')
$c.NumberFormat = '@'
$c.Value = 'BTSEToken'
$c.Style = 'Normal'
$c = $ws.Range('C20')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '2.320'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '12.90%'
$c.Style = 'Normal'
$c = $ws.Range('G20')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('B21')
$c.NumberFormat = '@'
$c.Value = 'BitpandaEcosystemToken'
$c.Style = 'Normal'
$c = $ws.Range('C21')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.3149'
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '1.17%'
$c.Style = 'Normal'
$c = $ws.Range('G21')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '-0.10%'
$c.Style = 'Normal'
$c = $ws.Range('G22')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '3.555'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '-4.95%'
$c.Style = 'Normal'
$c = $ws.Range('G23')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.1374'
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '0.00%'
$c.Style = 'Normal'
$c = $ws.Range('G24')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.001215'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '-2.65%'
$c.Style = 'Normal'
$c = $ws.Range('G25')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.004480'
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '-2.01%'
$c.Style = 'Normal'
$c = $ws.Range('G26')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '22.95%'
$c.Style = 'Normal'
$c = $ws.Range('G27')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '-27.43%'
$c.Style = 'Normal'
$c = $ws.Range('G28')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G29')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G30')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G31')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G32')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G33')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G34')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G35')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G36')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G37')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G38')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G39')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '1.58%'
$c.Style = 'Normal'
$c = $ws.Range('G40')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.005970'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '-3.21%'
$c.Style = 'Normal'
$c = $ws.Range('G41')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.1059'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '0.28%'
$c.Style = 'Normal'
$c = $ws.Range('G42')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.002376'
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '-9.16%'
$c.Style = 'Normal'
$c = $ws.Range('G43')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.009709'
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '17.58%'
$c.Style = 'Normal'
$c = $ws.Range('G44')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.00005097'
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '-5.58%'
$c.Style = 'Normal'
$c = $ws.Range('G45')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '0.03%'
$c.Style = 'Normal'
$c = $ws.Range('G46')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '-30.34%'
$c.Style = 'Normal'
$c = $ws.Range('G47')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.002512'
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '3.41%'
$c.Style = 'Normal'
$c = $ws.Range('G48')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.00002100'
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '0.03%'
$c.Style = 'Normal'
$c = $ws.Range('G49')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0002000'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '0.03%'
$c.Style = 'Normal'
$c = $ws.Range('G50')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
$c = $ws.Range('G51')
$c.NumberFormat = '@'
$c.Value = '22'
$c.Style = 'Normal'
